# Apply cryptos list price/volume updates (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.592.80"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "2.278.87"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "123.85"
$ws.Range("E5").Value = "  +7.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "266.63"
$ws.Range("E6").Value = "  -1.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.641"
$ws.Range("E7").Value = "  +2.26%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.624"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.18"
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0947"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.30"
$ws.Range("E12").Value = "  +2.81%  "
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.47"
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.908"
$ws.Range("E15").Value = "  +3.32%  "
$ws.Range("D16").Value = "2.624.19"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").Value = "2.284.46"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").Value = "43.621.58"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.99"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.93"
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("E24").Value = "  -5.33%  "
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.97"
$ws.Range("E26").Value = "  +1.96%  "
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.73"
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.36"
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.92"
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.73"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0920"
$ws.Range("E33").Value = "  -1.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.78"
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.130"
$ws.Range("E35").Value = "  +1.91%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.29"
$ws.Range("E36").Value = "  +10.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0378"
$ws.Range("E37").Value = "  +4.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.65"
$ws.Range("E38").Value = "  -3.00%  "
$ws.Range("E39").Value = "  -1.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.55"
$ws.Range("E40").Value = "  +4.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.97"
$ws.Range("E41").Value = "  -3.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.46"
$ws.Range("E42").Value = "  -1.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.239"
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("E45").Value = "  -2.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.63"
$ws.Range("E46").Value = "  -11.88%  "
$ws.Range("B47").Value = "ordi"
$ws.Range("C47").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "73.80"
$ws.Range("E47").Value = "  +36.26%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.62"
$ws.Range("E48").Value = "  -1.64%  "
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "102.22"
$ws.Range("E51").Value = "  -0.79%  "
